$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.673.64"
$ws.Range("E2").Value = "  +3.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.210.54"
$ws.Range("E3").Value = "  +2.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.44"
$ws.Range("E5").Value = "  +4.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.39"
$ws.Range("E6").Value = "  +4.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.203.33"
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  +5.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +2.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.518"
$ws.Range("E12").Value = "  +4.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  +4.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.16"
$ws.Range("E14").Value = "  +6.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.721.11"
$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.372.15"
$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.44"
$ws.Range("E17").Value = "  +5.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.193.33"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.21"
$ws.Range("E20").Value = "  +2.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.33"
$ws.Range("E21").Value = "  +3.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("E22").Value = "  +4.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.10"
$ws.Range("E23").Value = "  +6.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.04"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.48"
$ws.Range("E25").Value = "  +2.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +5.35%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.02"
$ws.Range("E28").Value = "  +5.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +5.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("E30").Value = "  +15.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("E31").Value = "  +5.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.34"
$ws.Range("E32").Value = "  +3.83%  "

$ws.Range("E33").Value = "  +4.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.18"
$ws.Range("E36").Value = "  +1.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "493.21"
$ws.Range("E37").Value = "  +5.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0896"
$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0422"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.91"
$ws.Range("E40").Value = "  +3.86%  "

$ws.Range("E41").Value = "  +6.23%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  -2.74%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.300"
$ws.Range("E43").Value = "  +7.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0666"
$ws.Range("E44").Value = "  +16.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.920.93"
$ws.Range("E45").Value = "  -3.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.52"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  +11.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  +4.94%  "
